$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Chips
$ws.Range("A2").Value = "Chips"
$ws.Range("B2").Value = 1359684
$ws.Range("C2").Value = "R"
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 20
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 10

# Row 3 - Chocolate
$ws.Range("A3").Value = "Chocolate"
$ws.Range("B3").Value = 55432
$ws.Range("C3").Value = "R"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 6
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 10
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 22

# Row 4 - Dark Chocolate
$ws.Range("A4").Value = "Dark Chocolate"
$ws.Range("B4").Value = 1324
$ws.Range("C4").Value = "R"
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 10
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 30

# Widen column A (bestFit) so the longer "Dark Chocolate" entry fits
$ws.Columns.Item(1).ColumnWidth = 13.8
